# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F15").Value = 320
$wsExhibit.Range("F19").Value = 211
$wsExhibit.Range("F24").Value = 1410
$wsExhibit.Range("F33").Value = 4318
$wsExhibit.Range("F34").Value = 4091
$wsExhibit.Range("F36").Value = 81
$wsExhibit.Range("F40").Value = 476
$wsExhibit.Range("F42").Value = 1315
$wsExhibit.Range("F47").Value = 63

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F17").Value = 320
$wsAll.Range("F24").Value = 1410
$wsAll.Range("F32").Value = 4318
$wsAll.Range("F38").Value = 476
$wsAll.Range("F43").Value = 1315
$wsAll.Range("F47").Value = 63
